$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plans")

# Remove spaces from the header row column names
$ws.Range("A1").Value = "PlanName"
$ws.Range("C1").Value = "LCOPrice"
$ws.Range("D1").Value = "BCPrice"
$ws.Range("E1").Value = "SDCount"
$ws.Range("F1").Value = "HDCount"
